$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8.144923048233769
$ws.Range("D2").Value = 5.985329318825532
$ws.Range("E2").Value = 10.89523397165546
$ws.Range("F2").Value = 34.9474195185609
$ws.Range("G2").Value = 48.65861782161312
$ws.Range("H2").Value = 18.13475001300898
$ws.Range("I2").Value = 29.88768468568118
$ws.Range("L2").Value = 8.786270916479364
$ws.Range("M2").Value = 26.45699879384756
$ws.Range("N2").Value = 17.1832245783244
$ws.Range("C3").Value = 8.1653814908217
$ws.Range("D3").Value = 6.012278217060886
$ws.Range("E3").Value = 10.94219729129977
$ws.Range("F3").Value = 34.47795914709903
$ws.Range("G3").Value = 47.68840119806242
$ws.Range("H3").Value = 18.06106498245533
$ws.Range("I3").Value = 29.67976392764365
$ws.Range("L3").Value = 8.82268776377561
$ws.Range("M3").Value = 25.67477183505134
$ws.Range("N3").Value = 16.9544645273653
$ws.Range("C4").Value = 8.178647712473959
$ws.Range("D4").Value = 6.030275488145848
$ws.Range("E4").Value = 10.97248175825003
$ws.Range("F4").Value = 34.20026634595728
$ws.Range("G4").Value = 47.10387656138299
$ws.Range("H4").Value = 18.02146801213613
$ws.Range("I4").Value = 29.56228479846685
$ws.Range("L4").Value = 8.84610930233535
$ws.Range("M4").Value = 25.18450569200137
$ws.Range("N4").Value = 16.81459059331498
$ws.Range("C5").Value = 8.184231483854369
$ws.Range("D5").Value = 6.037972228221902
$ws.Range("E5").Value = 10.98518848811181
$ws.Range("F5").Value = 34.08989428093353
$ws.Range("G5").Value = 46.86885652797798
$ws.Range("H5").Value = 18.00675798569495
$ws.Range("I5").Value = 29.51700585018686
$ws.Range("L5").Value = 8.85592167801482
$ws.Range("M5").Value = 24.98249217192863
$ws.Range("N5").Value = 16.75780171870771
$ws.Range("C6").Value = 8.185169410471884
$ws.Range("D6").Value = 6.039272113470514
$ws.Range("E6").Value = 10.98732054675809
$ws.Range("F6").Value = 34.07173939699412
$ws.Range("G6").Value = 46.83003415834271
$ws.Range("H6").Value = 18.00440167086971
$ws.Range("I6").Value = 29.509644907444
$ws.Range("L6").Value = 8.85756722579548
$ws.Range("M6").Value = 24.94882217883395
$ws.Range("N6").Value = 16.74838656899079
$ws.Range("C7").Value = 8.178722297102023
$ws.Range("D7").Value = 6.030377823154174
$ws.Range("E7").Value = 10.97265164382559
$ws.Range("F7").Value = 34.19876635985146
$ws.Range("G7").Value = 47.10069365127787
$ws.Range("H7").Value = 18.02126384667199
$ws.Range("I7").Value = 29.56166360422
$ws.Range("L7").Value = 8.846240549379287
$ws.Range("M7").Value = 25.18178989236918
$ws.Range("N7").Value = 16.81382378061651
$ws.Range("C8").Value = 8.151831162240878
$ws.Range("D8").Value = 5.994319000943004
$ws.Range("E8").Value = 10.91112698249997
$ws.Range("F8").Value = 34.78344158510624
$ws.Range("G8").Value = 48.32198358138392
$ws.Range("H8").Value = 18.10817604474682
$ws.Range("I8").Value = 29.81390152275779
$ws.Range("L8").Value = 8.798607746551534
$ws.Range("M8").Value = 26.18952045894661
$ws.Range("N8").Value = 17.10426464141373
$ws.Range("C9").Value = 8.104666551647407
$ws.Range("D9").Value = 5.935210177142753
$ws.Range("E9").Value = 10.80191531658053
$ws.Range("F9").Value = 36.00749666671455
$ws.Range("G9").Value = 50.78966589712039
$ws.Range("H9").Value = 18.32303280332421
$ws.Range("I9").Value = 30.38762056348939
$ws.Range("L9").Value = 8.713575942356767
$ws.Range("M9").Value = 28.07516988456292
$ws.Range("N9").Value = 17.67569215225874
$ws.Range("C10").Value = 8.073378300717064
$ws.Range("D10").Value = 5.898986026059093
$ws.Range("E10").Value = 10.72856918284329
$ws.Range("F10").Value = 36.94534218860485
$ws.Range("G10").Value = 52.62573985771135
$ws.Range("H10").Value = 18.50732982982465
$ws.Range("I10").Value = 30.85468333748777
$ws.Range("L10").Value = 8.656143933959621
$ws.Range("M10").Value = 29.39166396672252
$ws.Range("N10").Value = 18.09300629732182
$ws.Range("C11").Value = 8.059868147973344
$ws.Range("D11").Value = 5.884099556066542
$ws.Range("E11").Value = 10.69668116406967
$ws.Range("F11").Value = 37.37845403399869
$ws.Range("G11").Value = 53.46168756791314
$ws.Range("H11").Value = 18.59673543996232
$ws.Range("I11").Value = 31.07639617290246
$ws.Range("L11").Value = 8.631097133926504
$ws.Range("M11").Value = 29.97313044034897
$ws.Range("N11").Value = 18.28155740915331
$ws.Range("C12").Value = 8.054855667372129
$ws.Range("D12").Value = 5.87869382340125
$ws.Range("E12").Value = 10.6848171362701
$ws.Range("F12").Value = 37.54323448998714
$ws.Range("G12").Value = 53.77799746736704
$ws.Range("H12").Value = 18.63137382122248
$ws.Range("I12").Value = 31.16162155182823
$ws.Range("L12").Value = 8.621766691925421
$ws.Range("M12").Value = 30.1906328115175
$ws.Range("N12").Value = 18.3527113444623
$ws.Range("C13").Value = 8.055930597767441
$ws.Range("D13").Value = 5.879847710381286
$ws.Range("E13").Value = 10.68736288943295
$ws.Range("F13").Value = 37.50771436044538
$ws.Range("G13").Value = 53.70989073391839
$ws.Range("H13").Value = 18.62387933625567
$ws.Range("I13").Value = 31.14321138086561
$ws.Range("L13").Value = 8.623769325158218
$ws.Range("M13").Value = 30.1439119623378
$ws.Range("N13").Value = 18.337398976157
$ws.Range("C14").Value = 8.059453696030106
$ws.Range("D14").Value = 5.883650170043657
$ws.Range("E14").Value = 10.69570087595302
$ws.Range("F14").Value = 37.39199610067016
$ws.Range("G14").Value = 53.48771742060446
$ws.Range("H14").Value = 18.59956959017394
$ws.Range("I14").Value = 31.08338267660719
$ws.Range("L14").Value = 8.630326427229621
$ws.Range("M14").Value = 29.99107930684202
$ws.Range("N14").Value = 18.28741656490861
$ws.Range("C15").Value = 8.061625163626426
$ws.Range("D15").Value = 5.886009503306338
$ws.Range("E15").Value = 10.7008356121141
$ws.Range("F15").Value = 37.32121088723989
$ws.Range("G15").Value = 53.35158787652104
$ws.Range("H15").Value = 18.58478049200626
$ws.Range("I15").Value = 31.04689905707186
$ws.Range("L15").Value = 8.634362899205863
$ws.Range("M15").Value = 29.8971099624605
$ws.Range("N15").Value = 18.25676704240331
$ws.Range("C16").Value = 8.074275731681428
$ws.Range("D16").Value = 5.899991142298077
$ws.Range("E16").Value = 10.73068276920611
$ws.Range("F16").Value = 36.91715388167954
$ws.Range("G16").Value = 52.57109290973083
$ws.Range("H16").Value = 18.50159748372553
$ws.Range("I16").Value = 30.84037468502266
$ws.Range("L16").Value = 8.657802437134707
$ws.Range("M16").Value = 29.35329733848007
$ws.Range("N16").Value = 18.08065285950781
$ws.Range("C17").Value = 8.08222130995912
$ws.Range("D17").Value = 5.908977967960308
$ws.Range("E17").Value = 10.74937060981408
$ws.Range("F17").Value = 36.67082258688598
$ws.Range("G17").Value = 52.09221381728011
$ws.Range("H17").Value = 18.45198105642971
$ws.Range("I17").Value = 30.71600329528067
$ws.Range("L17").Value = 8.672457583999611
$ws.Range("M17").Value = 29.01509006775362
$ws.Range("N17").Value = 17.97223864745011
$ws.Range("C18").Value = 8.086859478168032
$ws.Range("D18").Value = 5.914296653519638
$ws.Range("E18").Value = 10.76025850230306
$ws.Range("F18").Value = 36.52975900557413
$ws.Range("G18").Value = 51.81686172333075
$ws.Range("H18").Value = 18.42396828935887
$ws.Range("I18").Value = 30.64534212252833
$ws.Range("L18").Value = 8.680988484237592
$ws.Range("M18").Value = 28.81893281464736
$ws.Range("N18").Value = 17.90976403940081
$ws.Range("C19").Value = 8.088441589062951
$ws.Range("D19").Value = 5.91612310450329
$ws.Range("E19").Value = 10.76396888948554
$ws.Range("F19").Value = 36.4821086756362
$ws.Range("G19").Value = 51.7236579365786
$ws.Range("H19").Value = 18.41457439448953
$ws.Range("I19").Value = 30.62156939570873
$ws.Range("L19").Value = 8.683894387917224
$ws.Range("M19").Value = 28.75224336430022
$ws.Range("N19").Value = 17.88859290729792
$ws.Range("C20").Value = 8.081368446653249
$ws.Range("D20").Value = 5.908005793474617
$ws.Range("E20").Value = 10.74736686536855
$ws.Range("F20").Value = 36.6969820164161
$ws.Range("G20").Value = 52.14318492988465
$ws.Range("H20").Value = 18.45720857113416
$ws.Range("I20").Value = 30.72915282755858
$ws.Range("L20").Value = 8.670887005480939
$ws.Range("M20").Value = 29.05126282513327
$ws.Range("N20").Value = 17.98379214855896
$ws.Range("C21").Value = 8.058416070711731
$ws.Range("D21").Value = 5.882526993742921
$ws.Range("E21").Value = 10.69324608278485
$ws.Range("F21").Value = 37.42596572926576
$ws.Range("G21").Value = 53.55298445171455
$ws.Range("H21").Value = 18.60668887109026
$ws.Range("I21").Value = 31.10092192293005
$ws.Range("L21").Value = 8.628396269133288
$ws.Range("M21").Value = 30.03604422841862
$ws.Range("N21").Value = 18.30210477165564
$ws.Range("C22").Value = 8.044018558831338
$ws.Range("D22").Value = 5.867225390289068
$ws.Range("E22").Value = 10.65910593187616
$ws.Range("F22").Value = 37.90681576224469
$ws.Range("G22").Value = 54.47279090413709
$ws.Range("H22").Value = 18.70893378844063
$ws.Range("I22").Value = 31.35125212168409
$ws.Range("L22").Value = 8.601524683309078
$ws.Range("M22").Value = 30.66392625446723
$ws.Range("N22").Value = 18.50867813617987
$ws.Range("C23").Value = 8.051647737744741
$ws.Range("D23").Value = 5.875267754057734
$ws.Range("E23").Value = 10.67721493733352
$ws.Range("F23").Value = 37.64982576166983
$ws.Range("G23").Value = 53.98212769155772
$ws.Range("H23").Value = 18.65395380198033
$ws.Range("I23").Value = 31.21699408059266
$ws.Range("L23").Value = 8.615784654436526
$ws.Range("M23").Value = 30.33030842640564
$ws.Range("N23").Value = 18.39857950089166
$ws.Range("C24").Value = 8.081753807619908
$ws.Range("D24").Value = 5.908444840179699
$ws.Range("E24").Value = 10.74827230960218
$ws.Range("F24").Value = 36.6851535955992
$ws.Range("G24").Value = 52.12014098209576
$ws.Range("H24").Value = 18.45484361442971
$ws.Range("I24").Value = 30.72320529798965
$ws.Range("L24").Value = 8.671596735468661
$ws.Range("M24").Value = 29.03491445676257
$ws.Range("N24").Value = 17.97856926028247
$ws.Range("C25").Value = 8.116832892016076
$ws.Range("D25").Value = 5.949945298271134
$ws.Range("E25").Value = 10.83024385772027
$ws.Range("F25").Value = 35.66894691913277
$ws.Range("G25").Value = 50.11652917486936
$ws.Range("H25").Value = 18.26021062517032
$ws.Range("I25").Value = 30.22420455622624
$ws.Range("L25").Value = 8.735689352051306
$ws.Range("M25").Value = 27.57623938656527
$ws.Range("N25").Value = 17.52128720317551
